$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.008.58"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "2.695.30"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("D9").Value = "2.719.21"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").Value = "3.188.82"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "60.771.74"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "2.691.90"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "0.0₃0818"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.55%  "
$ws.Range("E36").Value = "  +9.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.951"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.61%  "
$ws.Range("E39").Value = "  +8.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "2.148.01"
$ws.Range("E47").Value = "  +7.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0541"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
